$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 259; this pushes the existing rows 259-282 down to 260-283
# and extends the used range from A1:T282 to A1:T283 (matching Excel's native
# "insert row" shifting behaviour for entire-row inserts).
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new weekly record. All of the
# "constant" columns for this Mango / Vega Modelo de Temuco sheet keep the
# same values as every other data row.
$ws.Cells.Item(259, 1).Value = 10
$ws.Cells.Item(259, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(259, 3).Value = "La Araucanía"
$ws.Cells.Item(259, 4).Value = 44578
$ws.Cells.Item(259, 5).Value = 9
$ws.Cells.Item(259, 6).Value = "Fruta"
$ws.Cells.Item(259, 7).Value = 100108
$ws.Cells.Item(259, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(259, 9).Value = 100108002
$ws.Cells.Item(259, 10).Value = "Mango"
$ws.Cells.Item(259, 11).Value = "Sin especificar"
$ws.Cells.Item(259, 12).Value = "Primera"
$ws.Cells.Item(259, 13).Value = 333
$ws.Cells.Item(259, 14).Value = 7000
$ws.Cells.Item(259, 15).Value = 8000
$ws.Cells.Item(259, 16).Value = 7646
$ws.Cells.Item(259, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(259, 18).Value = "Perú"
$ws.Cells.Item(259, 19).Value = 1912
$ws.Cells.Item(259, 20).Value = 4
